# Apply the "ProposedNames" Bengali (bn) translation refresh:
#  - Update the bn_new (column H) proposed-name values for a batch of rows
#  - Add a new blank "Sheet2" after "Sheet1" and make it the active sheet
#  - Leave Sheet1 with an all-cells selection and a slightly wider column E

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update the proposed Bengali (bn_new) names in column H ---
$bnUpdates = @{
    "H3"   = "অশোক"
    "H10"  = "অশোক"
    "H21"  = "বিবেক"
    "H25"  = "পার্থ"
    "H34"  = "ধীরজ"
    "H43"  = "অলিভিয়া"
    "H45"  = "সইকত"
    "H47"  = "বিনোদ "
    "H62"  = "মিষ্টি "
    "H63"  = "মনোজ "
    "H64"  = "আশীষ "
    "H68"  = "অবির"
    "H78"  = "অভীক "
    "H84"  = "নিলয় "
    "H100" = "মৌলি "
    "H113" = "অর্ক"
    "H124" = "নীলাভ্র"
    "H128" = "দেবেন "
}

foreach ($addr in $bnUpdates.Keys) {
    $ws1.Range($addr).Value = $bnUpdates[$addr]
}

# --- Cosmetic view-state tweaks observed on Sheet1 ---
# (Target stored column width is 10.85546875 characters; ColumnWidth maps to the
# stored <col width> with a fixed ~5/6-character padding offset in this engine, so
# 10.0 is the input that lands closest to that stored value.)
$ws1.Columns("E:E").ColumnWidth = 10
$ws1.Cells.Select()

# --- Add a new empty Sheet2 right after Sheet1 and make it active ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"
$ws2.Range("E9").Select()
